$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 4).Value = -7.989000000000002
$ws.Cells.Item(4, 1).Value = -20.872
$ws.Cells.Item(6, 1).Value = -22.184
$ws.Cells.Item(6, 5).Value = 16.373
$ws.Cells.Item(7, 1).Value = -21.59500000000001
$ws.Cells.Item(7, 2).Value = 6.355
$ws.Cells.Item(8, 1).Value = -21.833
$ws.Cells.Item(8, 5).Value = 16.308
$ws.Cells.Item(11, 2).Value = 6.709000000000001
$ws.Cells.Item(12, 2).Value = 5.56
$ws.Cells.Item(12, 3).Value = -12.068
$ws.Cells.Item(12, 4).Value = -7.631
$ws.Cells.Item(13, 3).Value = -12.825
$ws.Cells.Item(13, 4).Value = -8.404999999999999
$ws.Cells.Item(14, 3).Value = -12.067
$ws.Cells.Item(15, 2).Value = 5.087000000000001
$ws.Cells.Item(16, 1).Value = -21.726
$ws.Cells.Item(16, 3).Value = -13.144
$ws.Cells.Item(18, 5).Value = 16.666
$ws.Cells.Item(19, 3).Value = -12.704
$ws.Cells.Item(19, 5).Value = 16.826
$ws.Cells.Item(20, 1).Value = -21.087
$ws.Cells.Item(20, 2).Value = 5.973999999999999
$ws.Cells.Item(20, 3).Value = -12.666
$ws.Cells.Item(20, 5).Value = 15.923
$ws.Cells.Item(21, 1).Value = -21.215
$ws.Cells.Item(21, 2).Value = 6.776999999999999
$ws.Cells.Item(21, 5).Value = 16.952
$ws.Cells.Item(22, 2).Value = 8.245000000000001
$ws.Cells.Item(22, 3).Value = -12.466
$ws.Cells.Item(22, 4).Value = -8.087
$ws.Cells.Item(23, 2).Value = 8.352
$ws.Cells.Item(23, 5).Value = 16.239
$ws.Cells.Item(24, 5).Value = 16.719
$ws.Cells.Item(25, 4).Value = -8.352
$ws.Cells.Item(28, 1).Value = -21.749
$ws.Cells.Item(29, 1).Value = -21.042
$ws.Cells.Item(29, 2).Value = 6.538000000000001
$ws.Cells.Item(29, 4).Value = -7.128
$ws.Cells.Item(30, 1).Value = -21.465
$ws.Cells.Item(32, 1).Value = -21.486
$ws.Cells.Item(34, 2).Value = 7.983
$ws.Cells.Item(34, 4).Value = -7.976999999999999
$ws.Cells.Item(35, 5).Value = 16.406
$ws.Cells.Item(36, 3).Value = -12.89
$ws.Cells.Item(37, 5).Value = 16.509
$ws.Cells.Item(39, 5).Value = 17.057
$ws.Cells.Item(40, 1).Value = -20.779
$ws.Cells.Item(41, 5).Value = 16.597
$ws.Cells.Item(42, 2).Value = 7.392
$ws.Cells.Item(43, 2).Value = 4.664
$ws.Cells.Item(43, 3).Value = -13.075
$ws.Cells.Item(43, 4).Value = -8.424000000000001
$ws.Cells.Item(44, 2).Value = 5.633
$ws.Cells.Item(45, 2).Value = 5.165999999999999
$ws.Cells.Item(46, 1).Value = -21.216
$ws.Cells.Item(46, 2).Value = 6.922
$ws.Cells.Item(46, 3).Value = -13.141
$ws.Cells.Item(46, 5).Value = 16.723
$ws.Cells.Item(47, 5).Value = 16.632
$ws.Cells.Item(48, 4).Value = -7.741000000000001
$ws.Cells.Item(48, 5).Value = 17.04
$ws.Cells.Item(50, 2).Value = 4.956
$ws.Cells.Item(50, 3).Value = -13.384
$ws.Cells.Item(51, 1).Value = -21.512
$ws.Cells.Item(51, 2).Value = 6.522
$ws.Cells.Item(52, 1).Value = -21.783
$ws.Cells.Item(57, 1).Value = -21.444
$ws.Cells.Item(57, 2).Value = 6.909000000000001
$ws.Cells.Item(57, 5).Value = 16.377
$ws.Cells.Item(58, 5).Value = 16.599
$ws.Cells.Item(59, 1).Value = -22.134
$ws.Cells.Item(60, 4).Value = -8.141000000000002
$ws.Cells.Item(60, 5).Value = 16.735
$ws.Cells.Item(62, 1).Value = -21.926
$ws.Cells.Item(65, 2).Value = 4.964
$ws.Cells.Item(66, 1).Value = -21.492
$ws.Cells.Item(66, 2).Value = 6.082000000000001
$ws.Cells.Item(67, 2).Value = 5.8
$ws.Cells.Item(68, 4).Value = -6.703999999999999
$ws.Cells.Item(70, 4).Value = -7.164999999999999
$ws.Cells.Item(71, 4).Value = -7.522999999999999
$ws.Cells.Item(73, 1).Value = -20.118
$ws.Cells.Item(73, 4).Value = -8.259
$ws.Cells.Item(73, 5).Value = 16.473
$ws.Cells.Item(74, 1).Value = -21.244
$ws.Cells.Item(76, 3).Value = -12.467
$ws.Cells.Item(76, 5).Value = 16.562
$ws.Cells.Item(77, 1).Value = -21.176
$ws.Cells.Item(78, 4).Value = -8.098000000000001
$ws.Cells.Item(79, 2).Value = 5.790999999999999
$ws.Cells.Item(84, 2).Value = 6.086
$ws.Cells.Item(85, 5).Value = 16.885
$ws.Cells.Item(87, 2).Value = 5.048
$ws.Cells.Item(87, 4).Value = -8.241000000000001
$ws.Cells.Item(92, 1).Value = -21.095
$ws.Cells.Item(92, 2).Value = 6.225
$ws.Cells.Item(92, 4).Value = -6.497
$ws.Cells.Item(95, 3).Value = -11.862
$ws.Cells.Item(97, 2).Value = 6.189
$ws.Cells.Item(97, 3).Value = -13.166
$ws.Cells.Item(98, 5).Value = 16.362
$ws.Cells.Item(99, 3).Value = -12.123
$ws.Cells.Item(100, 1).Value = -21.526
$ws.Cells.Item(101, 4).Value = -7.867
$ws.Cells.Item(101, 5).Value = 16.495
